$d = $word.ActiveDocument

# --- Change 1: "{% for semester in semesters %}" -> "{% for semester in courses %}" ---
$rng1 = $d.Content
$rng1.Find.Execute("semesters", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1.Text = "courses"

# --- Change 2: "{% elif semester.ZET_check %}" -> "{% elif semester.ZET_check == '2' %}" ---
$rng2 = $d.Content
$rng2.Find.Execute("ZET_check %}", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
# "ZET_check %}" -> insertion point right after "ZET_check " (10 chars: Z E T _ c h e c k <space>)
$insertAt = $rng2.Start + 10
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertAfter("== '2' ")

# Force the newly inserted text into its own run (distinct from the neighboring,
# now-identical-looking, runs) by nudging a character property away and back.
$newRng = $d.Range($insertAt, $insertAt + 7)
$newRng.Font.Size = 99
$newRng.Font.Size = 14
